# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The "Periodo Mora" column (E16:E33) lists the overdue billing periods for
# the worker, and the matching "Valor Mora" column (F16:F33) lists the debt
# amount for that period. The old statement periods are being retired and
# replaced: the period list is reversed so the most recent period (2104)
# now appears first (row 16) and the oldest (1911) now appears last (row 33).
# The one period whose arrears amount differs from the rest (28708, on the
# last row) moves along with its original period, so after the reversal it
# still accompanies the "1911" period (now on the last row again).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original order of periods as they appeared in rows 16-33 (top to bottom).
$periods = @(
    "1911","1912","2001","2002","2003","2004","2005","2006","2007",
    "2008","2009","2010","2011","2012","2101","2102","2103","2104"
)

# Matching "Valor Mora" amounts for those same rows (top to bottom).
$valores = @(
    33125,33125,33125,33125,33125,33125,33125,33125,33125,
    33125,33125,33125,33125,33125,33125,33125,33125,28708
)

$firstRow = 16
$lastRow = 33
$rowCount = $lastRow - $firstRow + 1

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $firstRow + $i
    $srcIndex = $rowCount - 1 - $i   # reverse the lists

    $ws.Cells.Item($row, 5).Value = $periods[$srcIndex]   # column E
    $ws.Cells.Item($row, 6).Value = $valores[$srcIndex]   # column F
}
